$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.016767382621765
$ws.Range("B1").Value = 1.678268313407898
$ws.Range("C1").Value = 4.36162281036377
$ws.Range("D1").Value = 2.417491436004639
$ws.Range("E1").Value = 1.356452226638794
